$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A631").Value = 540425
$ws.Range("B631").Value = "2025-09-12T18:30:00Z"
$ws.Range("C631").Value = 3
$ws.Range("D631").Value = 19
$ws.Range("E631").Value = "Bayer 04 Leverkusen"
$ws.Range("F631").Value = 19
$ws.Range("G631").Value = "Eintracht Frankfurt"
$ws.Range("H631").Value = 3
$ws.Range("I631").Value = 1
$ws.Range("J631").Value = "HomeWin"

$ws.Range("A632").Value = 540426
$ws.Range("B632").Value = "2025-09-13T13:30:00Z"
$ws.Range("C632").Value = 3
$ws.Range("D632").Value = 17
$ws.Range("E632").Value = "SC Freiburg"
$ws.Range("F632").Value = 10
$ws.Range("G632").Value = "VfB Stuttgart"
$ws.Range("H632").Value = 3
$ws.Range("I632").Value = 1
$ws.Range("J632").Value = "HomeWin"

$ws.Range("A633").Value = 540427
$ws.Range("B633").Value = "2025-09-13T13:30:00Z"
$ws.Range("C633").Value = 3
$ws.Range("D633").Value = 15
$ws.Range("E633").Value = "1. FSV Mainz 05"
$ws.Range("F633").Value = 721
$ws.Range("G633").Value = "RB Leipzig"
$ws.Range("H633").Value = 0
$ws.Range("I633").Value = 1
$ws.Range("J633").Value = "AwayWin"

$ws.Range("A634").Value = 540429
$ws.Range("B634").Value = "2025-09-13T13:30:00Z"
$ws.Range("C634").Value = 3
$ws.Range("D634").Value = 11
$ws.Range("E634").Value = "VfL Wolfsburg"
$ws.Range("F634").Value = 1
$ws.Range("G634").Value = "1. FC Köln"
$ws.Range("H634").Value = 3
$ws.Range("I634").Value = 3
$ws.Range("J634").Value = "Draw"

$ws.Range("A635").Value = 540430
$ws.Range("B635").Value = "2025-09-13T13:30:00Z"
$ws.Range("C635").Value = 3
$ws.Range("D635").Value = 28
$ws.Range("E635").Value = "1. FC Union Berlin"
$ws.Range("F635").Value = 2
$ws.Range("G635").Value = "TSG 1899 Hoffenheim"
$ws.Range("H635").Value = 2
$ws.Range("I635").Value = 4
$ws.Range("J635").Value = "AwayWin"

$ws.Range("A636").Value = 540432
$ws.Range("B636").Value = "2025-09-13T13:30:00Z"
$ws.Range("C636").Value = 3
$ws.Range("D636").Value = 44
$ws.Range("E636").Value = "1. FC Heidenheim 1846"
$ws.Range("F636").Value = 4
$ws.Range("G636").Value = "Borussia Dortmund"
$ws.Range("H636").Value = 0
$ws.Range("I636").Value = 2
$ws.Range("J636").Value = "AwayWin"

$ws.Range("A637").Value = 540424
$ws.Range("B637").Value = "2025-09-13T16:30:00Z"
$ws.Range("C637").Value = 3
$ws.Range("D637").Value = 5
$ws.Range("E637").Value = "FC Bayern München"
$ws.Range("F637").Value = 7
$ws.Range("G637").Value = "Hamburger SV"
$ws.Range("H637").Value = 5
$ws.Range("I637").Value = 0
$ws.Range("J637").Value = "HomeWin"

$ws.Range("A638").Value = 540431
$ws.Range("B638").Value = "2025-09-14T13:30:00Z"
$ws.Range("C638").Value = 3
$ws.Range("D638").Value = 20
$ws.Range("E638").Value = "FC St. Pauli 1910"
$ws.Range("F638").Value = 16
$ws.Range("G638").Value = "FC Augsburg"
$ws.Range("H638").Value = 2
$ws.Range("I638").Value = 1
$ws.Range("J638").Value = "HomeWin"

$ws.Range("A639").Value = 540428
$ws.Range("B639").Value = "2025-09-14T15:30:00Z"
$ws.Range("C639").Value = 3
$ws.Range("D639").Value = 18
$ws.Range("E639").Value = "Borussia Mönchengladbach"
$ws.Range("F639").Value = 12
$ws.Range("G639").Value = "SV Werder Bremen"
$ws.Range("H639").Value = 0
$ws.Range("I639").Value = 4
$ws.Range("J639").Value = "AwayWin"
